$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.171.41'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.87%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.173.35'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.14'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.606'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '70.19'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -4.90%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.579'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -6.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.07'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -8.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0928'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.23%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.48%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -5.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.497.47'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.93'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.25%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.177.71'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '40.970.27'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.20%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -7.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.45'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.69%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.06'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '225.82'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.74%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -7.64%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.87'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -5.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.51'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.27%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.35%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.88'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.93'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.11'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +6.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0768'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.17'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -9.03%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.95%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -8.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.13'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.32%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -5.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.41'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.51%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.26%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '60.30'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -7.66%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.33'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0971'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '98.30'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.90%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.52%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.24'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -7.10%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.374.47'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.16%  '
